# Final pass on Test_Scenario.xlsx: update the Keyword_Scenario sheet's
# Module_Reference values for the last two rows (Sr No 2 & 3) and make
# Keyword_Scenario the active/selected sheet again (it had drifted to
# Test_Scenarios), leaving the cursor on C4.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Keyword_Scenario")

$ws.Range("C3").Value = "Application_Submit1"
$ws.Range("C4").Value = "Application_Submit1"

$ws.Activate()
$ws.Range("C4").Select()
